$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header labels: "..._old" -> "..._FV2310", "..._new" -> "..._FV2404"
#    Columns A:J hold the "_old" headers, K holds "diff" (unchanged),
#    L:U hold the "_new" headers.
# ---------------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $colOld = $i + 1        # A..J
    $colNew = $i + 12       # L..U
    $ws.Cells.Item(1, $colOld).Value = ($oldHeaders[$i] + "_FV2310")
    $ws.Cells.Item(1, $colNew).Value = ($oldHeaders[$i] + "_FV2404")
}

# ---------------------------------------------------------------------------
# 2. Turn A1:U79 into an Excel table ("Table1") without disturbing the
#    existing header-row formatting (bold/fill/border already baked into the
#    sheet's cell style) and without Excel minting a new header dxf for it.
#    We stash the current header formatting, reset the header to the default
#    style so table creation has nothing "custom" to capture, build the
#    table, then restore the original formatting verbatim.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A200:U200")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)  # xlPasteFormats
$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$scratchRange.Clear()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1, top-left cell of the scrolling
#    pane is A2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
